$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "F6"  = 100
    "F8"  = 400
    "F9"  = 300
    "F10" = 900
    "F11" = 400
    "F13" = 500
    "F15" = 900
    "F16" = 500
    "F17" = 500
    "F18" = 600
    "F19" = 550
    "F20" = 550
    "F22" = 900
    "F23" = 600
    "F24" = 900
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
